$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with changed values.
# Price cells are set to Text number format first so numeric-looking strings
# (e.g. "1.00", "0.350") keep their exact formatting instead of being coerced
# into numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.197.75"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.643.41"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.06"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.90"
$ws.Range("E6").Value = "  +3.43%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.542"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.141"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.350"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.92"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.125.39"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000186"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.191.70"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.641.69"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.32"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "358.80"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.40"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.72"
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.75"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.73"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.813.40"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "554.75"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.37"
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.86"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -3.77%  "
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.370"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.32"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.60"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0318"
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "156.47"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.78"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.91"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.611"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.567"
$ws.Range("E51").Value = "  -0.45%  "

# Rows 37 and 38 swap content: Monero moves up to row 37 (was EthereumClassic),
# EthereumClassic moves to row 38 (was Monero), each with updated price/volume.
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.61"
$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.66"
$ws.Range("E38").Value = "  +1.31%  "
